$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.575.62"
$ws.Range("E2").Value = "  -0.61%  "
Set-TextValue "D3" "1.884.97"
$ws.Range("E3").Value = "  -0.57%  "
Set-TextValue "D4" "1.004"
$ws.Range("E4").Value = "  +0.63%  "
Set-TextValue "D5" "235.79"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("E6").Value = "  +0.68%  "
Set-TextValue "D7" "0.4856"
$ws.Range("E7").Value = "  -1.48%  "
Set-TextValue "D8" "0.2886"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("E9").Value = "  -1.15%  "
Set-TextValue "D10" "1.889.67"
$ws.Range("E10").Value = "  -0.25%  "
Set-TextValue "D11" "16.72"
$ws.Range("E11").Value = "  -0.01%  "
Set-TextValue "D12" "0.07206"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +0.45%  "
Set-TextValue "D14" "4.979"
$ws.Range("E14").Value = "  -1.27%  "
Set-TextValue "D15" "0.6629"
Set-TextValue "D16" "30.533.35"
$ws.Range("E16").Value = "  -0.23%  "
Set-TextValue "D17" "0.000007819"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("E18").Value = "  +0.33%  "
Set-TextValue "D19" "12.96"
Set-TextValue "D20" "2.133.88"
$ws.Range("E20").Value = "  +0.34%  "
Set-TextValue "D21" "1.005"
$ws.Range("E21").Value = "  +0.57%  "
Set-TextValue "D22" "4.743"
$ws.Range("E22").Value = "  -1.22%  "
Set-TextValue "D23" "185.48"
$ws.Range("E23").Value = "  +18.51%  "
Set-TextValue "D24" "5.996"
$ws.Range("E24").Value = "  +0.97%  "
Set-TextValue "D25" "9.225"
$ws.Range("E25").Value = "  -0.37%  "
Set-TextValue "D26" "155.12"
$ws.Range("E26").Value = "  +2.18%  "
Set-TextValue "D27" "18.44"
$ws.Range("E27").Value = "  +7.24%  "
Set-TextValue "D28" "1.848"
$ws.Range("E28").Value = "  -3.92%  "
Set-TextValue "D29" "1.406"
$ws.Range("E29").Value = "  -0.16%  "
Set-TextValue "D30" "4.226"
$ws.Range("E30").Value = "  -0.83%  "
Set-TextValue "D31" "0.08991"
$ws.Range("E31").Value = "  +1.73%  "
Set-TextValue "D32" "3.900"
Set-TextValue "D33" "0.05229"
$ws.Range("E33").Value = "  -0.10%  "
Set-TextValue "D34" "0.7269"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "1.076"
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.721"
$ws.Range("E36").Value = "  +2.16%  "
Set-TextValue "D37" "0.01810"
$ws.Range("E37").Value = "  -1.76%  "
Set-TextValue "D38" "2.666"
$ws.Range("E38").Value = "  -1.20%  "
Set-TextValue "D39" "0.9148"
$ws.Range("E39").Value = "  -2.70%  "
Set-TextValue "D40" "2.058"
$ws.Range("E40").Value = "  -6.01%  "
Set-TextValue "D41" "0.4310"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D42" "104.08"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D43" "0.9997"
$ws.Range("E43").Value = "  +0.03%  "
Set-TextValue "D44" "5.596"
$ws.Range("E44").Value = "  -4.26%  "
Set-TextValue "D45" "0.1326"
$ws.Range("E45").Value = "  +1.58%  "
Set-TextValue "D46" "7.310"
$ws.Range("E46").Value = "  -2.99%  "
Set-TextValue "D47" "0.05843"
$ws.Range("E47").Value = "  +0.79%  "
Set-TextValue "D48" "8.686"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "1.401"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D50" "33.25"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D51" "0.3877"
$ws.Range("E51").Value = "  +1.79%  "
